# regenerate orders with updated distance/sizes
#
# The experiment's distance conditions and one of the sizes were renumbered:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31   (S25 / S20 unchanged)
#
# These codes appear embedded inside several text columns (Condition,
# Filename_Left, Filename_Right, Distance, Size), so walk every used cell
# and rewrite any string value that contains one of the old tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne $null) {
            if ($v.GetType().Name -eq "String") {
                $nv = $v.Replace("D64", "D69").Replace("D80", "D86").Replace("D51", "D55").Replace("S30", "S31")
                if ($nv -ne $v) {
                    $cell.Value2 = $nv
                }
            }
        }
    }
}
